# scratchpad.xlsx edit
#
# The "Who Knows About Who" data-stores table gets a new leading column
# (A) that classifies each data store as "static" / "struct" / "dynamic".
# The previous columns A:D shuffle right into B:E; columns F:G (the
# "Movement" / id-key notes) are untouched. Two new rows (10, 11) are
# appended describing a new "IComponent"/"Mobility" dynamic store.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 title: was merged A1:B1 ("Who Knows About Who" + blank, centered
# style). Now it's merged B1:C1 instead; A1 goes back to a plain empty
# default cell.
# ---------------------------------------------------------------------
$ws.Range("A1:B1").UnMerge()

$ws.Range("B1").Value2 = $ws.Range("A1").Value2
$ws.Range("B1").HorizontalAlignment = -4108   # xlCenter

$ws.Range("A1").Value2 = ""
$ws.Range("A1").Style = "Normal"

$ws.Range("C1").Value2 = ""
$ws.Range("C1").HorizontalAlignment = -4108   # xlCenter

$ws.Range("B1:C1").Merge()

# ---------------------------------------------------------------------
# Data rows: write the full new grid directly (columns A:E), cell by
# cell, rather than trying to "shift" the old contents, since several
# rows move to different rows too (A3/A5/A6/A8 all collapse to a
# repeated "static"/"Map" marker etc).
# ---------------------------------------------------------------------

# Row 3
$ws.Range("A3").Value2 = "static"
$ws.Range("B3").Value2 = "Map"
$ws.Range("C3").Value2 = "No one"

# Row 4
$ws.Range("A4").Value2 = "struct"
$ws.Range("B4").Value2 = "XY"
$ws.Range("C4").Value2 = "No one"

# Row 5
$ws.Range("A5").Value2 = "static"
$ws.Range("B5").Value2 = "EntityManager"
$ws.Range("C5").Value2 = "Entity"
$ws.Range("D5").Value2 = "Component"

# Row 6
$ws.Range("A6").Value2 = "static"
$ws.Range("B6").Value2 = "MovementController"
$ws.Range("C6").Value2 = "EntityManager"
$ws.Range("D6").Value2 = "Component"
$ws.Range("E6").Value2 = "Map"

# Row 8
$ws.Range("A8").Value2 = "static"
$ws.Range("B8").Value2 = "DisplayProcessor"
$ws.Range("C8").Value2 = "EntityManager"
$ws.Range("D8").Value2 = "Movement"
$ws.Range("E8").Value2 = "Map"

# Row 9
$ws.Range("A9").Value2 = "dynamic"
$ws.Range("B9").Value2 = "Entity"
$ws.Range("C9").Value2 = "No one"

# Row 10 (new)
$ws.Range("A10").Value2 = "dynamic"
$ws.Range("B10").Value2 = "IComponent"
$ws.Range("C10").Value2 = "No one"

# Row 11 (new)
$ws.Range("A11").Value2 = "dynamic"
$ws.Range("B11").Value2 = "Mobility"
$ws.Range("C11").Value2 = "Movement"

# Columns F/G ("Who Knows About Who" side of the sheet) keep their
# values/positions untouched - nothing to do there.

# ---------------------------------------------------------------------
# Column widths: old A/B/C widths (20.421875 / 13.421875 / 11.00390625)
# now belong to the new B/C/D columns; column A reverts to the sheet's
# standard width.
# ---------------------------------------------------------------------
$ws.Range("B:B").ColumnWidth = 19.588541666666668
$ws.Range("C:C").ColumnWidth = 12.588541666666666
$ws.Range("D:D").ColumnWidth = 10.170572916666666
$ws.Range("A:A").ColumnWidth = $ws.StandardWidth

Write-Output "scratchpad edit applied"
